$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row above the current row 1. This pushes the header row
# (formerly row 1) down to row 2 and the existing data rows down by one as well.
$ws.Rows.Item(1).Insert()

# Turn the freshly-inserted row 1 into a title/banner row: centered text,
# a bit taller than the default row, merged across the full A:G width.
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(1, $c).HorizontalAlignment = -4108
}
$ws.Rows.Item(1).RowHeight = 27.75
$ws.Range("A1:G1").Merge()

# Move the active selection to D7 (matches the post-edit cursor position).
$ws.Range("D7").Select()
